$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '03/19/2021'
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("Q2").Value = '$/caja 16 kilos'
$ws.Range("S2").Value = 875
$ws.Range("T2").Value = 16
# Row 3
$ws.Range("D3").Value = '03/19/2021'
$ws.Range("M3").Value = 130
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = '$/caja 16 kilos'
$ws.Range("S3").Value = 750
$ws.Range("T3").Value = 16
# Row 5
$ws.Range("D5").Value = '03/12/2021'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("S5").Value = 722
# Row 6
$ws.Range("D6").Value = '01/28/2021'
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = '$/caja 16 kilos'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 1125
$ws.Range("T6").Value = 16
# Row 7
$ws.Range("D7").Value = '04/12/2021'
$ws.Range("L7").Value = 'Extra (doble especial)'
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 1111
# Row 8
$ws.Range("D8").Value = '03/18/2021'
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 16 kilos'
$ws.Range("S8").Value = 938
$ws.Range("T8").Value = 16
# Row 9
$ws.Range("D9").Value = '03/18/2021'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("Q9").Value = '$/caja 16 kilos'
$ws.Range("S9").Value = 812
$ws.Range("T9").Value = 16
# Row 10
$ws.Range("D10").Value = '03/18/2021'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("R10").Value = 'Provincia de Melipilla'
$ws.Range("S10").Value = 625
# Row 11
$ws.Range("D11").Value = '03/22/2021'
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 833
# Row 12
$ws.Range("D12").Value = '02/25/2021'
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("S12").Value = 722
# Row 13
$ws.Range("D13").Value = '04/05/2021'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("S13").Value = 1000
# Row 14
$ws.Range("D14").Value = '04/06/2021'
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("S14").Value = 889
# Row 15
$ws.Range("D15").Value = '04/06/2021'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("S15").Value = 778
# Row 16
$ws.Range("D16").Value = '04/14/2021'
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("S16").Value = 1000
# Row 17
$ws.Range("D17").Value = '04/14/2021'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("S17").Value = 889
$ws.Range("T17").Value = 18
# Row 18
$ws.Range("D18").Value = '04/09/2021'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("S18").Value = 556
$ws.Range("T18").Value = 18
# Row 19
$ws.Range("D19").Value = '04/13/2021'
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 170
$ws.Range("Q19").Value = '$/caja 18 kilos'
$ws.Range("R19").Value = 'Provincia de Melipilla'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18
# Row 20
$ws.Range("D20").Value = '04/13/2021'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 16000
$ws.Range("S20").Value = 889
# Row 21
$ws.Range("D21").Value = '01/26/2021'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 18000
$ws.Range("Q21").Value = '$/caja 16 kilos'
$ws.Range("S21").Value = 1125
$ws.Range("T21").Value = 16
# Row 22
$ws.Range("D22").Value = '04/15/2021'
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 16000
$ws.Range("O22").Value = 16000
$ws.Range("P22").Value = 16000
$ws.Range("S22").Value = 889
# Row 23
$ws.Range("D23").Value = '03/29/2021'
$ws.Range("L23").Value = 'Especial'
$ws.Range("M23").Value = 120
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("S23").Value = 722
# Row 24
$ws.Range("D24").Value = '03/29/2021'
$ws.Range("L24").Value = 'Extra (doble especial)'
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 15000
$ws.Range("P24").Value = 15000
$ws.Range("S24").Value = 833
# Row 25
$ws.Range("D25").Value = '03/29/2021'
$ws.Range("M25").Value = 50
# Row 26
$ws.Range("D26").Value = '03/24/2021'
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 14000
$ws.Range("O26").Value = 14000
$ws.Range("P26").Value = 14000
$ws.Range("R26").Value = 'Provincia de Melipilla'
$ws.Range("S26").Value = 778
# Row 27
$ws.Range("D27").Value = '03/24/2021'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 12000
$ws.Range("O27").Value = 12000
$ws.Range("P27").Value = 12000
$ws.Range("S27").Value = 667
# Row 28
$ws.Range("D28").Value = '06/03/2021'
$ws.Range("L28").Value = 'Especial'
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("R28").Value = 'Provincia de Limarí'
$ws.Range("S28").Value = 1333
# Row 29
$ws.Range("D29").Value = '04/29/2021'
$ws.Range("L29").Value = 'Especial'
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 24000
$ws.Range("O29").Value = 24000
$ws.Range("P29").Value = 24000
$ws.Range("S29").Value = 1333
# Row 30
$ws.Range("D30").Value = '04/29/2021'
$ws.Range("N30").Value = 20000
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 20000
$ws.Range("S30").Value = 1111
# Row 31
$ws.Range("D31").Value = '05/14/2021'
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 23000
$ws.Range("O31").Value = 23000
$ws.Range("P31").Value = 23000
$ws.Range("S31").Value = 1278
# Row 32
$ws.Range("D32").Value = '03/16/2021'
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 60
$ws.Range("Q32").Value = '$/caja 18 kilos'
$ws.Range("S32").Value = 833
$ws.Range("T32").Value = 18
# Row 33
$ws.Range("D33").Value = '04/23/2021'
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 20000
$ws.Range("Q33").Value = '$/caja 18 kilos'
$ws.Range("S33").Value = 1111
$ws.Range("T33").Value = 18
# Row 34
$ws.Range("D34").Value = '04/23/2021'
$ws.Range("L34").Value = 'Primera'
$ws.Range("N34").Value = 18000
$ws.Range("O34").Value = 18000
$ws.Range("P34").Value = 18000
$ws.Range("Q34").Value = '$/caja 18 kilos'
$ws.Range("S34").Value = 1000
$ws.Range("T34").Value = 18
